$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "43.046.04"
Set-TextValue "E2" "  +0.15%  "
Set-TextValue "D3" "2.309.17"
Set-TextValue "E3" "  +0.07%  "
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "301.60"
Set-TextValue "E5" "  -0.82%  "
Set-TextValue "D6" "98.50"
Set-TextValue "E6" "  -1.69%  "
Set-TextValue "E7" "  +0.70%  "
Set-TextValue "E8" "  +0.05%  "
Set-TextValue "D9" "0.522"
Set-TextValue "E9" "  +1.66%  "
Set-TextValue "D10" "35.65"
Set-TextValue "E10" "  +2.02%  "
Set-TextValue "D11" "0.0788"
Set-TextValue "E11" "  -1.00%  "
Set-TextValue "E12" "  -1.11%  "
Set-TextValue "D13" "17.90"
Set-TextValue "E13" "  -0.53%  "
Set-TextValue "D14" "6.88"
Set-TextValue "E14" "  +0.02%  "
Set-TextValue "D15" "2.669.71"
Set-TextValue "E15" "  -0.55%  "
Set-TextValue "D16" "2.302.95"
Set-TextValue "E16" "  +0.00%  "
Set-TextValue "E17" "  -3.50%  "
Set-TextValue "D18" "42.979.41"
Set-TextValue "E18" "  +0.14%  "
Set-TextValue "D19" "13.36"
Set-TextValue "E19" "  +7.13%  "
Set-TextValue "D20" "6.18"
Set-TextValue "E20" "  +0.50%  "
Set-TextValue "E21" "  +0.51%  "
Set-TextValue "D22" "68.00"
Set-TextValue "E22" "  +0.45%  "
Set-TextValue "D23" "240.77"
Set-TextValue "E23" "  +1.39%  "
Set-TextValue "E24" "  -2.73%  "
Set-TextValue "B25" "Dai"
Set-TextValue "C25" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D25" "0.999"
Set-TextValue "E25" "  -0.06%  "
Set-TextValue "B26" "PancakeSwap"
Set-TextValue "C26" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D26" "2.44"
Set-TextValue "E26" "  -1.25%  "
Set-TextValue "D27" "24.97"
Set-TextValue "E27" "  +0.70%  "
Set-TextValue "D28" "169.01"
Set-TextValue "E28" "  +0.72%  "
Set-TextValue "B29" "Toncoin"
Set-TextValue "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D29" "2.05"
Set-TextValue "E29" "  -10.48%  "
Set-TextValue "B30" "Cosmos"
Set-TextValue "C30" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D30" "9.16"
Set-TextValue "E30" "  -0.12%  "
Set-TextValue "D31" "33.36"
Set-TextValue "E31" "  -1.81%  "
Set-TextValue "E32" "  +3.81%  "
Set-TextValue "D33" "4.88"
Set-TextValue "E33" "  +5.45%  "
Set-TextValue "E34" "  -0.05%  "
Set-TextValue "D35" "18.29"
Set-TextValue "E35" "  +7.51%  "
Set-TextValue "E36" "  -0.07%  "
Set-TextValue "E37" "  +0.38%  "
Set-TextValue "E38" "  +0.37%  "
Set-TextValue "E39" "  +0.86%  "
Set-TextValue "D40" "2.75"
Set-TextValue "E40" "  -2.54%  "
Set-TextValue "E41" "  -0.48%  "
Set-TextValue "D42" "1.989.49"
Set-TextValue "E42" "  -0.69%  "
Set-TextValue "E43" "  +1.15%  "
Set-TextValue "D44" "10.11"
Set-TextValue "E44" "  -0.98%  "
Set-TextValue "D45" "2.07"
Set-TextValue "E45" "  -10.29%  "
Set-TextValue "D46" "17.47"
Set-TextValue "E46" "  +0.71%  "
Set-TextValue "E47" "  -0.46%  "
Set-TextValue "D48" "75.87"
Set-TextValue "E48" "  +8.31%  "
Set-TextValue "D49" "54.55"
Set-TextValue "E49" "  -1.55%  "
Set-TextValue "D50" "2.537.51"
Set-TextValue "E50" "  +0.65%  "
Set-TextValue "E51" "  -0.24%  "
